$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Cells.Item(1, 1).Value = 'Date Posted'
$ws.Cells.Item(1, 2).Value = 'Title'
$ws.Cells.Item(1, 3).Value = 'URL'
$ws.Cells.Item(1, 4).Value = 'Score'
$ws.Cells.Item(1, 5).Value = 'Date Scraped'

# Copy header style (bold, border, centered) from A1 onto the new D1/E1 header cells
$ws.Range("A1").Copy() | Out-Null
$ws.Range("D1:E1").PasteSpecial(-4122) | Out-Null

# --- Data rows ---
# Row 2
$ws.Cells.Item(2, 1).NumberFormat = '@'
$ws.Cells.Item(2, 1).Value = '2025-03-07'
$ws.Cells.Item(2, 2).Value = 'PREMIUM LITE? REALLY YOUTUBE?'
$ws.Cells.Item(2, 3).Value = 'https://i.redd.it/ps8fhv549cne1.jpeg'
$ws.Cells.Item(2, 4).Value = 7619
$ws.Cells.Item(2, 5).NumberFormat = '@'
$ws.Cells.Item(2, 5).Value = '2025-03-09'
$ws.Range("B2").Copy() | Out-Null
$ws.Range("A2").PasteSpecial(-4122) | Out-Null
$ws.Range("E2").PasteSpecial(-4122) | Out-Null

# Row 3
$ws.Cells.Item(3, 1).NumberFormat = '@'
$ws.Cells.Item(3, 1).Value = '2025-03-04'
$ws.Cells.Item(3, 2).Value = 'I Thought it was Supposed to be Bigger?'
$ws.Cells.Item(3, 3).Value = 'https://i.redd.it/2zukjeomekme1.jpeg'
$ws.Cells.Item(3, 4).Value = 7410
$ws.Cells.Item(3, 5).NumberFormat = '@'
$ws.Cells.Item(3, 5).Value = '2025-03-09'
$ws.Range("B3").Copy() | Out-Null
$ws.Range("A3").PasteSpecial(-4122) | Out-Null
$ws.Range("E3").PasteSpecial(-4122) | Out-Null

# Row 4
$ws.Cells.Item(4, 1).NumberFormat = '@'
$ws.Cells.Item(4, 1).Value = '2025-03-06'
$ws.Cells.Item(4, 2).Value = 'How can someone donate $500 in low effort Shorts just like that?'
$ws.Cells.Item(4, 3).Value = 'https://i.redd.it/g1gwn7bbozme1.jpeg'
$ws.Cells.Item(4, 4).Value = 5854
$ws.Cells.Item(4, 5).NumberFormat = '@'
$ws.Cells.Item(4, 5).Value = '2025-03-09'
$ws.Range("B4").Copy() | Out-Null
$ws.Range("A4").PasteSpecial(-4122) | Out-Null
$ws.Range("E4").PasteSpecial(-4122) | Out-Null

# Row 5
$ws.Cells.Item(5, 1).NumberFormat = '@'
$ws.Cells.Item(5, 1).Value = '2025-03-03'
$ws.Cells.Item(5, 2).Value = 'Ublock Origin is gone.'
$ws.Cells.Item(5, 3).Value = 'https://www.reddit.com/r/youtube/comments/1j2ec76/ublock_origin_is_gone/'
$ws.Cells.Item(5, 4).Value = 4609
$ws.Cells.Item(5, 5).NumberFormat = '@'
$ws.Cells.Item(5, 5).Value = '2025-03-09'
$ws.Range("B5").Copy() | Out-Null
$ws.Range("A5").PasteSpecial(-4122) | Out-Null
$ws.Range("E5").PasteSpecial(-4122) | Out-Null

# Row 6
$ws.Cells.Item(6, 1).NumberFormat = '@'
$ws.Cells.Item(6, 1).Value = '2025-03-03'
$ws.Cells.Item(6, 2).Value = 'Stay in School'
$ws.Cells.Item(6, 3).Value = 'https://www.reddit.com/gallery/1j2ve1c'
$ws.Cells.Item(6, 4).Value = 2635
$ws.Cells.Item(6, 5).NumberFormat = '@'
$ws.Cells.Item(6, 5).Value = '2025-03-09'
$ws.Range("B6").Copy() | Out-Null
$ws.Range("A6").PasteSpecial(-4122) | Out-Null
$ws.Range("E6").PasteSpecial(-4122) | Out-Null

# Row 7
$ws.Cells.Item(7, 1).NumberFormat = '@'
$ws.Cells.Item(7, 1).Value = '2025-03-02'
$ws.Cells.Item(7, 2).Value = 'The fact that he doesn''t even address the costume in the video is hilarious'
$ws.Cells.Item(7, 3).Value = 'https://i.redd.it/ytrmm5srbbme1.jpeg'
$ws.Cells.Item(7, 4).Value = 2286
$ws.Cells.Item(7, 5).NumberFormat = '@'
$ws.Cells.Item(7, 5).Value = '2025-03-09'
$ws.Range("B7").Copy() | Out-Null
$ws.Range("A7").PasteSpecial(-4122) | Out-Null
$ws.Range("E7").PasteSpecial(-4122) | Out-Null

# Row 8
$ws.Cells.Item(8, 1).NumberFormat = '@'
$ws.Cells.Item(8, 1).Value = '2025-03-04'
$ws.Cells.Item(8, 2).Value = 'I painted "Dorito Wars" by VanossGaming'
$ws.Cells.Item(8, 3).Value = 'https://i.redd.it/j1b3r48i7pme1.jpeg'
$ws.Cells.Item(8, 4).Value = 2114
$ws.Cells.Item(8, 5).NumberFormat = '@'
$ws.Cells.Item(8, 5).Value = '2025-03-09'
$ws.Range("B8").Copy() | Out-Null
$ws.Range("A8").PasteSpecial(-4122) | Out-Null
$ws.Range("E8").PasteSpecial(-4122) | Out-Null

# Row 9
$ws.Cells.Item(9, 1).NumberFormat = '@'
$ws.Cells.Item(9, 1).Value = '2025-03-02'
$ws.Cells.Item(9, 2).Value = 'Who is the YouTuber you would rate the lowest?'
$ws.Cells.Item(9, 3).Value = 'https://i.redd.it/cuucddozc9me1.png'
$ws.Cells.Item(9, 4).Value = 1918
$ws.Cells.Item(9, 5).NumberFormat = '@'
$ws.Cells.Item(9, 5).Value = '2025-03-09'
$ws.Range("B9").Copy() | Out-Null
$ws.Range("A9").PasteSpecial(-4122) | Out-Null
$ws.Range("E9").PasteSpecial(-4122) | Out-Null

# Row 10
$ws.Cells.Item(10, 1).NumberFormat = '@'
$ws.Cells.Item(10, 1).Value = '2025-03-06'
$ws.Cells.Item(10, 2).Value = 'I think my channel got shadow banned'
$ws.Cells.Item(10, 3).Value = 'https://i.redd.it/5j2k17r0b5ne1.jpeg'
$ws.Cells.Item(10, 4).Value = 1833
$ws.Cells.Item(10, 5).NumberFormat = '@'
$ws.Cells.Item(10, 5).Value = '2025-03-09'
$ws.Range("B10").Copy() | Out-Null
$ws.Range("A10").PasteSpecial(-4122) | Out-Null
$ws.Range("E10").PasteSpecial(-4122) | Out-Null

# Row 11
$ws.Cells.Item(11, 1).NumberFormat = '@'
$ws.Cells.Item(11, 1).Value = '2025-03-06'
$ws.Cells.Item(11, 2).Value = 'Bro waited 10 years for a song title'
$ws.Cells.Item(11, 3).Value = 'https://i.redd.it/ourctbrop5ne1.jpeg'
$ws.Cells.Item(11, 4).Value = 1785
$ws.Cells.Item(11, 5).NumberFormat = '@'
$ws.Cells.Item(11, 5).Value = '2025-03-09'
$ws.Range("B11").Copy() | Out-Null
$ws.Range("A11").PasteSpecial(-4122) | Out-Null
$ws.Range("E11").PasteSpecial(-4122) | Out-Null

# Row 12
$ws.Cells.Item(12, 1).NumberFormat = '@'
$ws.Cells.Item(12, 1).Value = '2025-03-05'
$ws.Cells.Item(12, 2).Value = 'Stop posting porn noone fucking wants to see it'
$ws.Cells.Item(12, 3).Value = 'https://www.reddit.com/r/youtube/comments/1j4egrw/stop_posting_porn_noone_fucking_wants_to_see_it/'
$ws.Cells.Item(12, 4).Value = 1291
$ws.Cells.Item(12, 5).NumberFormat = '@'
$ws.Cells.Item(12, 5).Value = '2025-03-09'
$ws.Range("B12").Copy() | Out-Null
$ws.Range("A12").PasteSpecial(-4122) | Out-Null
$ws.Range("E12").PasteSpecial(-4122) | Out-Null

# Row 13
$ws.Cells.Item(13, 1).NumberFormat = '@'
$ws.Cells.Item(13, 1).Value = '2025-03-08'
$ws.Cells.Item(13, 2).Value = 'Popular YouTuber who doesn''t deserve their fame?'
$ws.Cells.Item(13, 3).Value = 'https://www.reddit.com/gallery/1j6mlcx'
$ws.Cells.Item(13, 4).Value = 1289
$ws.Cells.Item(13, 5).NumberFormat = '@'
$ws.Cells.Item(13, 5).Value = '2025-03-09'
$ws.Range("B13").Copy() | Out-Null
$ws.Range("A13").PasteSpecial(-4122) | Out-Null
$ws.Range("E13").PasteSpecial(-4122) | Out-Null

# Row 14
$ws.Cells.Item(14, 1).NumberFormat = '@'
$ws.Cells.Item(14, 1).Value = '2025-03-08'
$ws.Cells.Item(14, 2).Value = 'Grandpa won’t stop watching AI-generated cheating stories at full volume on his TV and he’s driving my Grandma insane.'
$ws.Cells.Item(14, 3).Value = 'https://www.reddit.com/r/youtube/comments/1j6o6u9/grandpa_wont_stop_watching_aigenerated_cheating/'
$ws.Cells.Item(14, 4).Value = 1290
$ws.Cells.Item(14, 5).NumberFormat = '@'
$ws.Cells.Item(14, 5).Value = '2025-03-09'
$ws.Range("B14").Copy() | Out-Null
$ws.Range("A14").PasteSpecial(-4122) | Out-Null
$ws.Range("E14").PasteSpecial(-4122) | Out-Null

# Row 15
$ws.Cells.Item(15, 1).NumberFormat = '@'
$ws.Cells.Item(15, 1).Value = '2025-03-05'
$ws.Cells.Item(15, 2).Value = 'Where''s the ''Nudity & Sexual content'' report button?'
$ws.Cells.Item(15, 3).Value = 'https://i.redd.it/dqfx2ln0rume1.jpeg'
$ws.Cells.Item(15, 4).Value = 1144
$ws.Cells.Item(15, 5).NumberFormat = '@'
$ws.Cells.Item(15, 5).Value = '2025-03-09'
$ws.Range("B15").Copy() | Out-Null
$ws.Range("A15").PasteSpecial(-4122) | Out-Null
$ws.Range("E15").PasteSpecial(-4122) | Out-Null

# Row 16
$ws.Cells.Item(16, 1).NumberFormat = '@'
$ws.Cells.Item(16, 1).Value = '2025-03-07'
$ws.Cells.Item(16, 2).Value = 'Who the FUCK asked for this'
$ws.Cells.Item(16, 3).Value = 'https://i.redd.it/pzqi0nwxlane1.png'
$ws.Cells.Item(16, 4).Value = 1138
$ws.Cells.Item(16, 5).NumberFormat = '@'
$ws.Cells.Item(16, 5).Value = '2025-03-09'
$ws.Range("B16").Copy() | Out-Null
$ws.Range("A16").PasteSpecial(-4122) | Out-Null
$ws.Range("E16").PasteSpecial(-4122) | Out-Null

# Row 17
$ws.Cells.Item(17, 1).NumberFormat = '@'
$ws.Cells.Item(17, 1).Value = '2025-03-02'
$ws.Cells.Item(17, 2).Value = 'This is the worst thumbnail i''ve ever seen'
$ws.Cells.Item(17, 3).Value = 'https://i.redd.it/sjdc5y21u8me1.png'
$ws.Cells.Item(17, 4).Value = 1076
$ws.Cells.Item(17, 5).NumberFormat = '@'
$ws.Cells.Item(17, 5).Value = '2025-03-09'
$ws.Range("B17").Copy() | Out-Null
$ws.Range("A17").PasteSpecial(-4122) | Out-Null
$ws.Range("E17").PasteSpecial(-4122) | Out-Null

# Row 18
$ws.Cells.Item(18, 1).NumberFormat = '@'
$ws.Cells.Item(18, 1).Value = '2025-03-06'
$ws.Cells.Item(18, 2).Value = 'Me when non-predatory Minecraft YouTuber'
$ws.Cells.Item(18, 3).Value = 'https://i.redd.it/vpeklffxfzme1.jpeg'
$ws.Cells.Item(18, 4).Value = 941
$ws.Cells.Item(18, 5).NumberFormat = '@'
$ws.Cells.Item(18, 5).Value = '2025-03-09'
$ws.Range("B18").Copy() | Out-Null
$ws.Range("A18").PasteSpecial(-4122) | Out-Null
$ws.Range("E18").PasteSpecial(-4122) | Out-Null

# Row 19
$ws.Cells.Item(19, 1).NumberFormat = '@'
$ws.Cells.Item(19, 1).Value = '2025-03-03'
$ws.Cells.Item(19, 2).Value = 'I just got a 43 minute unskippable add what the hell youtube'
$ws.Cells.Item(19, 3).Value = 'https://i.redd.it/0axg0m2oqeme1.jpeg'
$ws.Cells.Item(19, 4).Value = 835
$ws.Cells.Item(19, 5).NumberFormat = '@'
$ws.Cells.Item(19, 5).Value = '2025-03-09'
$ws.Range("B19").Copy() | Out-Null
$ws.Range("A19").PasteSpecial(-4122) | Out-Null
$ws.Range("E19").PasteSpecial(-4122) | Out-Null

# Row 20
$ws.Cells.Item(20, 1).NumberFormat = '@'
$ws.Cells.Item(20, 1).Value = '2025-03-02'
$ws.Cells.Item(20, 2).Value = 'People are flooding Markiplier''s video about him losing his niece in a car crash with disgusting comments. These aren''t even a quarter of the ones made on that video.'
$ws.Cells.Item(20, 3).Value = 'https://www.reddit.com/gallery/1j1tdiy'
$ws.Cells.Item(20, 4).Value = 827
$ws.Cells.Item(20, 5).NumberFormat = '@'
$ws.Cells.Item(20, 5).Value = '2025-03-09'
$ws.Range("B20").Copy() | Out-Null
$ws.Range("A20").PasteSpecial(-4122) | Out-Null
$ws.Range("E20").PasteSpecial(-4122) | Out-Null

# Row 21
$ws.Cells.Item(21, 1).NumberFormat = '@'
$ws.Cells.Item(21, 1).Value = '2025-03-07'
$ws.Cells.Item(21, 2).Value = 'People who use AI voice to narrate their videos are so annoying.'
$ws.Cells.Item(21, 3).Value = 'https://www.reddit.com/r/youtube/comments/1j5us1j/people_who_use_ai_voice_to_narrate_their_videos/'
$ws.Cells.Item(21, 4).Value = 613
$ws.Cells.Item(21, 5).NumberFormat = '@'
$ws.Cells.Item(21, 5).Value = '2025-03-09'
$ws.Range("B21").Copy() | Out-Null
$ws.Range("A21").PasteSpecial(-4122) | Out-Null
$ws.Range("E21").PasteSpecial(-4122) | Out-Null

Write-Output "Applied reddit_posts edit: headers + 20 data rows (A1:E21)"
